# Auto-generated update of commodity price/profit columns (H-N) across multiple crafting sheets
# Source: scheduled market-data refresh

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 56600
$ws.Range("I21").Value = 44900
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 44900
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -44432
$ws.Range("N21").Value = -80936

$ws.Range("H23").Value = 56600
$ws.Range("I23").Value = 44900
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 44900
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -44666
$ws.Range("N23").Value = -80468

$ws.Range("H28").Value = 699.6667
$ws.Range("I28").Value = 632.6923
$ws.Range("J28").Value = 873.8
$ws.Range("K28").Value = 632.6923
$ws.Range("L28").Value = 873.8
$ws.Range("M28").Value = -147.6923
$ws.Range("N28").Value = -1843.8

$ws.Range("H32").Value = 973.5
$ws.Range("I32").Value = 710.25
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 710.25
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -384.25
$ws.Range("N32").Value = -2152

$ws.Range("H39").Value = 1031.4445
$ws.Range("I39").Value = 145.75
$ws.Range("J39").Value = 1740
$ws.Range("K39").Value = 437.25
$ws.Range("L39").Value = 5220
$ws.Range("M39").Value = -141.25
$ws.Range("N39").Value = -5812

$ws.Range("H43").Value = 1800.2354
$ws.Range("I43").Value = 795
$ws.Range("J43").Value = 1934.2667
$ws.Range("K43").Value = 795
$ws.Range("L43").Value = 1934.2667
$ws.Range("M43").Value = -726
$ws.Range("N43").Value = -2072.2667

$ws.Range("H51").Value = 6408.273
$ws.Range("I51").Value = 3800.3333
$ws.Range("J51").Value = 7386.25
$ws.Range("K51").Value = 3800.3333
$ws.Range("L51").Value = 7386.25
$ws.Range("M51").Value = -3316.3333
$ws.Range("N51").Value = -8354.25

$ws.Range("H81").Value = 40325.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 40325.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 40325.5
$ws.Range("N81").Value = -42321.5

$ws.Range("H84").Value = 40325.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 40325.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 120976.5
$ws.Range("N84").Value = -130960.5

$ws.Range("H93").Value = 30000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H98").Value = 2316.6667
$ws.Range("I98").Value = 2380
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2380
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -882
$ws.Range("N98").Value = -4996

$ws.Range("H101").Value = 562.0909
$ws.Range("I101").Value = 400.57144
$ws.Range("J101").Value = 844.75
$ws.Range("K101").Value = 1201.71432
$ws.Range("L101").Value = 2534.25
$ws.Range("M101").Value = 420.28568
$ws.Range("N101").Value = -5778.25

$ws.Range("H105").Value = 34950
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 34950
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 34950
$ws.Range("N105").Value = -41938

$ws.Range("H106").Value = 5186.7856
$ws.Range("I106").Value = 6326.875
$ws.Range("J106").Value = 3666.6667
$ws.Range("K106").Value = 6326.875
$ws.Range("L106").Value = 3666.6667
$ws.Range("M106").Value = -5695.875
$ws.Range("N106").Value = -4928.6667

$ws.Range("H107").Value = 1656.7333
$ws.Range("I107").Value = 338.14285
$ws.Range("J107").Value = 2810.5
$ws.Range("K107").Value = 338.14285
$ws.Range("L107").Value = 2810.5
$ws.Range("M107").Value = 1581.85715
$ws.Range("N107").Value = -6650.5

$ws.Range("H122").Value = 2316.6667
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4690
$ws.Range("N122").Value = -10900

$ws.Range("H125").Value = 1511.1428
$ws.Range("I125").Value = 741.3333
$ws.Range("J125").Value = 1721.091
$ws.Range("K125").Value = 6671.9997
$ws.Range("L125").Value = 15489.819
$ws.Range("M125").Value = -4211.9997
$ws.Range("N125").Value = -20409.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20436.518
$ws.Range("I32").Value = 18568.537
$ws.Range("J32").Value = 29643
$ws.Range("K32").Value = 18568.537
$ws.Range("L32").Value = 29643
$ws.Range("M32").Value = -18281.537
$ws.Range("N32").Value = -30217

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1474.7142
$ws.Range("I64").Value = 696.6667
$ws.Range("J64").Value = 2058.25
$ws.Range("K64").Value = 696.6667
$ws.Range("L64").Value = 2058.25
$ws.Range("M64").Value = -471.6667
$ws.Range("N64").Value = -2508.25

$ws.Range("H67").Value = 1474.7142
$ws.Range("I67").Value = 696.6667
$ws.Range("J67").Value = 2058.25
$ws.Range("K67").Value = 696.6667
$ws.Range("L67").Value = 2058.25
$ws.Range("M67").Value = 83.33330000000001
$ws.Range("N67").Value = -3618.25

$ws.Range("H107").Value = 1500.2812
$ws.Range("I107").Value = 1384.2174
$ws.Range("J107").Value = 1796.8889
$ws.Range("K107").Value = 1384.2174
$ws.Range("L107").Value = 1796.8889
$ws.Range("M107").Value = 535.7826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1233.3334
$ws.Range("I105").Value = 1180
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1180
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 567

$ws.Range("H107").Value = 596.1177
$ws.Range("I107").Value = 531.38464
$ws.Range("J107").Value = 806.5
$ws.Range("K107").Value = 531.38464
$ws.Range("L107").Value = 806.5
$ws.Range("M107").Value = 1388.61536

$ws.Range("H132").Value = 31262.03
$ws.Range("I132").Value = 1613.2333
$ws.Range("J132").Value = 253628
$ws.Range("K132").Value = 4839.699900000001
$ws.Range("L132").Value = 760884
$ws.Range("M132").Value = -2309.699900000001
$ws.Range("N132").Value = -765944

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5924.2856
$ws.Range("I3").Value = 5303.077
$ws.Range("J3").Value = 14000
$ws.Range("K3").Value = 15909.231
$ws.Range("L3").Value = 42000
$ws.Range("M3").Value = -15797.231
$ws.Range("N3").Value = -42224

$ws.Range("H12").Value = 35714372
$ws.Range("I12").Value = 100000080
$ws.Range("J12").Value = 89.44444
$ws.Range("K12").Value = 300000240
$ws.Range("L12").Value = 268.33332
$ws.Range("M12").Value = -300000067
$ws.Range("N12").Value = -614.33332

$ws.Range("H98").Value = 540
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 766.6667
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 2300.0001
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -5296.0001

$ws.Range("H141").Value = 10645.111
$ws.Range("I141").Value = 5110.1816
$ws.Range("J141").Value = 19342.857
$ws.Range("K141").Value = 15330.5448
$ws.Range("L141").Value = 58028.571
$ws.Range("M141").Value = -10150.5448
$ws.Range("N141").Value = -68388.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4244.3335
$ws.Range("I80").Value = 4050
$ws.Range("J80").Value = 4283.2
$ws.Range("K80").Value = 4050
$ws.Range("L80").Value = 4283.2
$ws.Range("M80").Value = -3052
$ws.Range("N80").Value = -6279.2

$ws.Range("H83").Value = 4244.3335
$ws.Range("I83").Value = 4050
$ws.Range("J83").Value = 4283.2
$ws.Range("K83").Value = 20250
$ws.Range("L83").Value = 21416
$ws.Range("M83").Value = -15258
$ws.Range("N83").Value = -31400

$ws.Range("H107").Value = 184.63637
$ws.Range("I107").Value = 158.66667
$ws.Range("J107").Value = 301.5
$ws.Range("K107").Value = 158.66667
$ws.Range("L107").Value = 301.5
$ws.Range("M107").Value = 1761.33333
$ws.Range("N107").Value = -4141.5

$ws.Range("H122").Value = 2086.4
$ws.Range("I122").Value = 1666.2858
$ws.Range("J122").Value = 3066.6667
$ws.Range("K122").Value = 4998.857400000001
$ws.Range("L122").Value = 9200.000100000001
$ws.Range("M122").Value = -2548.857400000001
$ws.Range("N122").Value = -14100.0001

$ws.Range("H126").Value = 1906.8422
$ws.Range("I126").Value = 2073.7778
$ws.Range("J126").Value = 1756.6
$ws.Range("K126").Value = 6221.3334
$ws.Range("L126").Value = 5269.799999999999
$ws.Range("M126").Value = -3751.3334
$ws.Range("N126").Value = -10209.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3290.158
$ws.Range("I40").Value = 3418.9375
$ws.Range("J40").Value = 2603.3333
$ws.Range("K40").Value = 3418.9375
$ws.Range("L40").Value = 2603.3333
$ws.Range("M40").Value = -3282.9375
$ws.Range("N40").Value = -2875.3333

$ws.Range("H82").Value = 2416.6924
$ws.Range("I82").Value = 1867.3334
$ws.Range("J82").Value = 2581.5
$ws.Range("K82").Value = 1867.3334
$ws.Range("L82").Value = 2581.5
$ws.Range("M82").Value = -1506.3334
$ws.Range("N82").Value = -3303.5

$ws.Range("H85").Value = 2416.6924
$ws.Range("I85").Value = 1867.3334
$ws.Range("J85").Value = 2581.5
$ws.Range("K85").Value = 1867.3334
$ws.Range("L85").Value = 2581.5
$ws.Range("M85").Value = -619.3334
$ws.Range("N85").Value = -5077.5

$ws.Range("H122").Value = 3431.3044
$ws.Range("I122").Value = 2930
$ws.Range("J122").Value = 3570.5557
$ws.Range("K122").Value = 8790
$ws.Range("L122").Value = 10711.6671
$ws.Range("M122").Value = -6340
$ws.Range("N122").Value = -15611.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 68283.664
$ws.Range("I132").Value = 51326.6
$ws.Range("J132").Value = 102197.8
$ws.Range("K132").Value = 153979.8
$ws.Range("L132").Value = 306593.4
$ws.Range("M132").Value = -151449.8
$ws.Range("N132").Value = -311653.4
